# Correct dates of research assistant entries on the "employment" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("employment")

# Row 2: Undergraduate Research Assistant (differential psychology) - start date
$ws.Range("E2").Value = "2020/10"

# Row 3: Undergraduate Research Assistant (educational psychology) - start + end date
$ws.Range("E3").Value = "2021/08"
$ws.Range("F3").Value = "2023/03"

# Row 4: Voluntary Social Year - reformat start/end date separators
$ws.Range("E4").Value = "2019/09"
$ws.Range("F4").Value = "2020/09"
